$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '44.005.52'
Set-TextCell $ws.Range("E2") '  +5.14%  '

# Row 3
Set-TextCell $ws.Range("D3") '2.276.80'
Set-TextCell $ws.Range("E3") '  +2.18%  '

# Row 4
Set-TextCell $ws.Range("E4") '  +0.14%  '

# Row 5
Set-TextCell $ws.Range("D5") '231.08'
Set-TextCell $ws.Range("E5") '  -0.22%  '

# Row 6
Set-TextCell $ws.Range("D6") '0.627'
Set-TextCell $ws.Range("E6") '  +0.64%  '

# Row 7
Set-TextCell $ws.Range("D7") '61.17'
Set-TextCell $ws.Range("E7") '  -0.65%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.423'
Set-TextCell $ws.Range("E9") '  +4.80%  '

# Row 10
Set-TextCell $ws.Range("D10") '57.99'
Set-TextCell $ws.Range("E10") '  -1.88%  '

# Row 11
Set-TextCell $ws.Range("D11") '0.0933'
Set-TextCell $ws.Range("E11") '  +4.05%  '

# Row 12
Set-TextCell $ws.Range("E12") '  +0.30%  '

# Row 13
Set-TextCell $ws.Range("D13") '2.614.31'
Set-TextCell $ws.Range("E13") '  +2.14%  '

# Row 14
Set-TextCell $ws.Range("D14") '23.69'
Set-TextCell $ws.Range("E14") '  +7.57%  '

# Row 15
Set-TextCell $ws.Range("D15") '15.53'
Set-TextCell $ws.Range("E15") '  -0.69%  '

# Row 16
Set-TextCell $ws.Range("D16") '5.82'
Set-TextCell $ws.Range("E16") '  +4.08%  '

# Row 17
Set-TextCell $ws.Range("D17") '0.808'
Set-TextCell $ws.Range("E17") '  +0.74%  '

# Row 18
Set-TextCell $ws.Range("D18") '2.274.56'
Set-TextCell $ws.Range("E18") '  +1.45%  '

# Row 19
Set-TextCell $ws.Range("D19") '43.867.87'
Set-TextCell $ws.Range("E19") '  +5.09%  '

# Row 20
Set-TextCell $ws.Range("E20") '  +3.59%  '

# Row 21
Set-TextCell $ws.Range("D21") '73.10'
Set-TextCell $ws.Range("E21") '  +1.43%  '

# Row 22
Set-TextCell $ws.Range("D22") '6.21'
Set-TextCell $ws.Range("E22") '  +2.98%  '

# Row 23
Set-TextCell $ws.Range("D23") '253.20'
Set-TextCell $ws.Range("E23") '  +1.56%  '

# Row 24
Set-TextCell $ws.Range("E24") '  +0.11%  '

# Row 25
Set-TextCell $ws.Range("D25") '2.55'
Set-TextCell $ws.Range("E25") '  +6.72%  '

# Row 26
Set-TextCell $ws.Range("E26") '  -0.75%  '

# Row 27
Set-TextCell $ws.Range("D27") '9.86'
Set-TextCell $ws.Range("E27") '  +2.65%  '

# Row 28
Set-TextCell $ws.Range("D28") '170.89'
Set-TextCell $ws.Range("E28") '  +1.80%  '

# Row 29
Set-TextCell $ws.Range("E29") '  -1.14%  '

# Row 30
Set-TextCell $ws.Range("D30") '20.50'
Set-TextCell $ws.Range("E30") '  +2.42%  '

# Row 31
Set-TextCell $ws.Range("E31") '  +2.13%  '

# Row 32
Set-TextCell $ws.Range("E32") '  -1.04%  '

# Row 33
Set-TextCell $ws.Range("E33") '  -0.81%  '

# Row 34
Set-TextCell $ws.Range("D34") '5.07'
Set-TextCell $ws.Range("E34") '  +0.25%  '

# Row 35
Set-TextCell $ws.Range("E35") '  +2.30%  '

# Row 36
Set-TextCell $ws.Range("E36") '  +3.11%  '

# Row 37
Set-TextCell $ws.Range("D37") '6.50'
Set-TextCell $ws.Range("E37") '  -1.86%  '

# Row 38
Set-TextCell $ws.Range("D38") '2.39'
Set-TextCell $ws.Range("E38") '  +0.96%  '

# Row 39
Set-TextCell $ws.Range("D39") '3.60'
Set-TextCell $ws.Range("E39") '  -2.20%  '

# Row 40
Set-TextCell $ws.Range("E40") '  +4.12%  '

# Row 41
Set-TextCell $ws.Range("E41") '  +0.49%  '

# Row 42
Set-TextCell $ws.Range("E42") '  +1.88%  '

# Row 43
Set-TextCell $ws.Range("E43") '  -14.91%  '

# Row 44
Set-TextCell $ws.Range("D44") '0.0985'
Set-TextCell $ws.Range("E44") '  +1.08%  '

# Row 45
Set-TextCell $ws.Range("D45") '4.49'
Set-TextCell $ws.Range("E45") '  -7.37%  '

# Row 46
Set-TextCell $ws.Range("D46") '98.41'
Set-TextCell $ws.Range("E46") '  -0.57%  '

# Row 47
Set-TextCell $ws.Range("E47") '  -1.18%  '

# Row 48
Set-TextCell $ws.Range("D48") '1.478.44'
Set-TextCell $ws.Range("E48") '  +0.04%  '

# Row 49
Set-TextCell $ws.Range("D49") '16.69'
Set-TextCell $ws.Range("E49") '  +1.35%  '

# Row 50
Set-TextCell $ws.Range("D50") '1.09'
Set-TextCell $ws.Range("E50") '  +0.67%  '

# Row 51
Set-TextCell $ws.Range("E51") '  -2.89%  '
